# [Outlook] (preview) Add mappings for new calendar properties (Compose)
# Adds 8 new rows to the "Snippets" table describing the new
# IsAllDayEvent / Sensitivity calendar property mappings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Data for the new rows, in order: Class, Member Name, Member ID (methods only), SnippetIdIntheYAMLFile, MethodNameInTheSnippet
$rows = @(
    @("AppointmentCompose", "isAllDayEvent", $null, "outlook-calendar-properties-apis", "getIsAllDayEvent"),
    @("IsAllDayEvent",      "getAsync",      2,     "outlook-calendar-properties-apis", "getIsAllDayEvent"),
    @("AppointmentCompose", "isAllDayEvent", $null, "outlook-calendar-properties-apis", "setIsAllDayEventTrue"),
    @("IsAllDayEvent",      "setAsync",      2,     "outlook-calendar-properties-apis", "setIsAllDayEventTrue"),
    @("AppointmentCompose", "sensitivity",   $null, "outlook-calendar-properties-apis", "getSensitivity"),
    @("Sensitivity",        "getAsync",      2,     "outlook-calendar-properties-apis", "getSensitivity"),
    @("AppointmentCompose", "sensitivity",   $null, "outlook-calendar-properties-apis", "setSensitivityConfidential"),
    @("Sensitivity",        "setAsync",      2,     "outlook-calendar-properties-apis", "setSensitivityConfidential")
)

foreach ($row in $rows) {
    $newRow = $lo.ListRows.Add()
    $r = $newRow.Range.Row

    # Make sure every cell in the row (even the blank "Member ID" one) is
    # materialized, matching the original author's fully-formatted rows.
    $ws.Rows($r).EntireRow.Style = "Normal"

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    if ($null -ne $row[2]) {
        $ws.Cells.Item($r, 3).Value = $row[2]
    }
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

# Move selection to the new last row, matching the author's saved view state.
$lastRow = $lo.Range.Row + $lo.Range.Rows.Count - 1
$ws.Cells.Item($lastRow, 1).Select()
